$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 2331
$ws.Range("I3").Value = 2484
$ws.Range("I4").Value = 609
$ws.Range("H5").Value = 801
$ws.Range("I5").Value = 222
$ws.Range("I6").Value = 2863
$ws.Range("H7").Value = 25971
$ws.Range("I7").Value = 8509

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 87
$ws.Range("I3").Value = 80
$ws.Range("I7").Value = 271

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 43
$ws.Range("I7").Value = 159

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 75
$ws.Range("I6").Value = 113
$ws.Range("I7").Value = 329

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 25
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I4").Value = 11
$ws.Range("I7").Value = 200

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I3").Value = 4
$ws.Range("I5").Value = 27
$ws.Range("I6").Value = 59
$ws.Range("I7").Value = 284
$ws.Range("I8").Value = 539
$ws.Range("I11").Value = 142
$ws.Range("I18").Value = 63
$ws.Range("I19").Value = 239
$ws.Range("I21").Value = 52
$ws.Range("I25").Value = 36
$ws.Range("I29").Value = 552
$ws.Range("I33").Value = 389
$ws.Range("I36").Value = 117
$ws.Range("I37").Value = 271
$ws.Range("I43").Value = 75
$ws.Range("I44").Value = 67
$ws.Range("I52").Value = 170
$ws.Range("I54").Value = 191
$ws.Range("I55").Value = 94
$ws.Range("I57").Value = 28
$ws.Range("I59").Value = 16
$ws.Range("I63").Value = 40
$ws.Range("I64").Value = 81
$ws.Range("I65").Value = 200
$ws.Range("I67").Value = 329
$ws.Range("I76").Value = 130
$ws.Range("I78").Value = 112
$ws.Range("I79").Value = 221
$ws.Range("I83").Value = 166
$ws.Range("I84").Value = 64
$ws.Range("I85").Value = 395
$ws.Range("H90").Value = 298
$ws.Range("I94").Value = 74
$ws.Range("I95").Value = 143
$ws.Range("I97").Value = 66
$ws.Range("I99").Value = 159
$ws.Range("I100").Value = 14
$ws.Range("H101").Value = 25971
$ws.Range("I101").Value = 8509

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 58
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 166

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 143

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 95
$ws.Range("I3").Value = 135
$ws.Range("I5").Value = 7
$ws.Range("I6").Value = 130
$ws.Range("I7").Value = 389

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 44
$ws.Range("I6").Value = 94
$ws.Range("I7").Value = 191

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 171
$ws.Range("I3").Value = 192
$ws.Range("I5").Value = 18
$ws.Range("I6").Value = 154
$ws.Range("I7").Value = 552

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 98
$ws.Range("I7").Value = 239

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I3").Value = 16
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I3").Value = 32
$ws.Range("I7").Value = 130

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 98
$ws.Range("I3").Value = 156
$ws.Range("I4").Value = 18
$ws.Range("I7").Value = 395

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 26
$ws.Range("I7").Value = 59

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I2").Value = 22
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 112

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value = 33
$ws.Range("I7").Value = 94

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I6").Value = 43
$ws.Range("I7").Value = 52

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I5").Value = 7
$ws.Range("I7").Value = 221

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 81

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I2").Value = 15
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 35
$ws.Range("I4").Value = 5
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 117

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("I5").Value = 7
$ws.Range("I6").Value = 14

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 43
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 170

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 74

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I2").Value = 9
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 60
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("I6").Value = 7
$ws.Range("I7").Value = 16

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I3").Value = 13
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 66

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 169
$ws.Range("I3").Value = 146
$ws.Range("I5").Value = 18
$ws.Range("I6").Value = 173
$ws.Range("I7").Value = 539

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 27

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("H5").Value = 12
$ws.Range("H7").Value = 298

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I2").Value = 9
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I6").Value = 43
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item('Andersonville')
$ws.Range("I2").Value = 1
$ws.Range("I6").Value = 4

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 86
$ws.Range("I7").Value = 284
